$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.950.15"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.888.75"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.07"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.75"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07995"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.75"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "1.895.39"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.001"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.181"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.018"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06776"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.11"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "27.970.10"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.99"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.362"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "2.120.63"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.52"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.496"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.53"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09569"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9645"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.648"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.362"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06122"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02247"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.213"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.229"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5975"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1907"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.35"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5695"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.945"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.398"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06864"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.41"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.070"
$ws.Range("E51").Value = "  -0.66%  "
